$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: strip the surrounding brackets from the error message and add the Statement ID (L2)
$ws.Range("K2").Value = "Please review the parentheses in the input statement. There is 1 additional opening parenthesis ('(')."
$ws.Range("L2").Value = "'1"

# Row 3: remove brackets from [OK], renumber statement id and logical linkage
$ws.Range("K3").Value = "OK"
$ws.Range("L3").Value = "'2.1"
$ws.Range("AO3").Value = "[OR].Bdir.[2.2]"

# Row 4: remove brackets from [OK], renumber statement id and logical linkage
$ws.Range("K4").Value = "OK"
$ws.Range("L4").Value = "'2.2"
$ws.Range("AO4").Value = "[OR].Bdir.[2.1]"

# Row 5: remove brackets from [OK], renumber statement id
$ws.Range("K5").Value = "OK"
$ws.Range("L5").Value = "'3"

# Row 6: remove brackets from [OK], renumber statement id and logical linkage
$ws.Range("K6").Value = "OK"
$ws.Range("L6").Value = "'4.1"
$ws.Range("AO6").Value = "[OR].Bdir.[4.2]"

# Row 7: remove brackets from [OK], renumber statement id and logical linkage
$ws.Range("K7").Value = "OK"
$ws.Range("L7").Value = "'4.2"
$ws.Range("AO7").Value = "[OR].Bdir.[4.1]"
